# Add 2022-Q1 sheet with fund holdings data, and update 总计 (totals) sheet
# with a new leading row summarizing 2022-Q1.

$wb = $excel.ActiveWorkbook

$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore  = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet handles appear to be position-based, so once the new
# sheet is inserted before it, the old "总计" handle now resolves to the
# newly inserted sheet instead. Re-fetch "总计" by name so subsequent
# edits land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header / column-A formatting from an existing quarter sheet so the
# new sheet matches the look (bold + bordered, centered) of its siblings.
$prevQuarter.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$prevQuarter.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows
$rows = @(
  @(0, "010963", "信达澳银周期动力混合",             "30.52", "89.82", "4.85", "1.4802", 4),
  @(1, "002350", "华安安华灵活配置混合",             "42.47", "93.61", "2.53", "1.0745", 5),
  @(2, "010363", "信达澳银匠心臻选两年持有期混合",     "50.40", "92.98", "1.52", "0.7661", 9),
  @(3, "014207", "华安产业精选混合A",                "27.31", "62.03", "1.99", "0.5435", 8),
  @(4, "014208", "华安产业精选混合C",                "7.93",  "62.03", "1.99", "0.1578", 8),
  @(5, "011160", "富国质量成长6个月持有期混合A",       "3.80",  "91.55", "2.26", "0.0859", 8),
  @(6, "011161", "富国质量成长6个月持有期混合C",       "0.12",  "91.55", "2.26", "0.0027", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = "'" + $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = "'" + $row[3]
    $newSheet.Range("E$r").Value = "'" + $row[4]
    $newSheet.Range("F$r").Value = "'" + $row[5]
    $newSheet.Range("G$r").Value = "'" + $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert a new leading data row into the "总计" (totals) sheet for
#    2022-Q1, pushing the existing rows down and renumbering column A.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 4.11

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Restore the originally active sheet/tab (adding a sheet shifts focus to
# it by default), so we don't introduce unrelated view-state changes.
$wb.Worksheets.Item("2021-Q3").Activate()
